$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'37.008.89"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.42%  "
$ws.Range("D3").Value = "'2.050.43"
$ws.Range("D3").Style = "Normal"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'245.88"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.91%  "
$ws.Range("D6").Value = "'0.658"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.69%  "
$ws.Range("D7").Value = "'58.58"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -4.46%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").Value = "'0.378"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.51%  "
$ws.Range("D10").Value = "'0.0773"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.71%  "
$ws.Range("E11").Value = "  +1.98%  "
$ws.Range("D12").Value = "'15.35"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.67%  "
$ws.Range("E13").Value = "  +7.74%  "
$ws.Range("D14").Value = "'2.351.97"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.31%  "
$ws.Range("D15").Value = "'5.73"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.31%  "
$ws.Range("D16").Value = "'2.071.48"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.59%  "
$ws.Range("D17").Value = "'18.15"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.15%  "
$ws.Range("D18").Value = "'36.978.52"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.46%  "
$ws.Range("D19").Value = "'73.93"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.33%  "
$ws.Range("D20").Value = "'0.0₃0886"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.95%  "
$ws.Range("D21").Value = "'5.41"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.99%  "
$ws.Range("D22").Value = "'237.88"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.56%  "
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("E24").Value = "  +1.14%  "
$ws.Range("E25").Value = "  +3.68%  "
$ws.Range("D26").Value = "'168.68"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.48%  "
$ws.Range("D27").Value = "'2.14"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.22%  "
$ws.Range("D28").Value = "'20.02"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.52%  "
$ws.Range("D29").Value = "'5.54"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +14.58%  "
$ws.Range("E30").Value = "  -1.20%  "
$ws.Range("D31").Value = "'1.12"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.40%  "
$ws.Range("D32").Value = "'4.69"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.12%  "
$ws.Range("E34").Value = "  +0.11%  "
$ws.Range("E35").Value = "  +6.02%  "
$ws.Range("D36").Value = "'0.0852"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.15%  "
$ws.Range("E37").Value = "  -1.22%  "
$ws.Range("E38").Value = "  -3.21%  "
$ws.Range("D39").Value = "'5.23"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.32%  "
$ws.Range("E40").Value = "  -3.17%  "
$ws.Range("B41").Value = "Cronos"
$ws.Range("C41").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D41").Value = "'0.0976"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -11.04%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").Value = "'0.0222"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.37%  "
$ws.Range("E43").Value = "  -0.02%  "
$ws.Range("D44").Value = "'97.76"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.38%  "
$ws.Range("D45").Value = "'16.95"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.83%  "
$ws.Range("D46").Value = "'1.299.80"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.27%  "
$ws.Range("D47").Value = "'2.37"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.48%  "
$ws.Range("E48").Value = "  -0.45%  "
$ws.Range("D49").Value = "'6.76"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.81%  "
$ws.Range("D50").Value = "'2.236.54"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.45%  "
$ws.Range("D51").Value = "'3.55"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.86%  "
